$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (A13 empty, B13/C13 = "8151869 - Livia Chaguri e Carvalho")
# is removed from the sheet; every row below it shifts up by one.
$ws.Rows.Item(13).Delete()

# After the shift, a handful of B/C cells hold content that differs from
# a plain shift (the source data was re-shuffled), so fix those explicitly.

# Row 10 (was "Objetivos:" / long PT objective text) -> now holds the
# "Docentes responsáveis" value.
$ws.Range("B10").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C10").Value = "8151869 - Livia Chaguri e Carvalho"

# Row 13 (now "Programa resumido:") -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (now "Programa:") -> "01/01/2018"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Row 18 (now "Método:") -> "Docentes responsáveis" value again
$ws.Range("B18").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C18").Value = "8151869 - Livia Chaguri e Carvalho"

# Row 19 (now "Critério:") -> "Aplicação de 2 provas (P1 e P2)."
$ws.Range("B19").Value = "Aplicação de 2 provas (P1 e P2)."
$ws.Range("C19").Value = "Aplicação de 2 provas (P1 e P2)."

# Row 20 (now "Norma de recuperação:") -> the MP grading-criteria paragraph
$criterioText = "A média do período (MP) será calculada por: MP = (P1+P2)/2. `nAlunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). `nAlunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). `nAlunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Row 21 (now "Bibliografia:") -> the "média final após recuperação" sentence
$normaText = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação"
$ws.Range("B21").Value = $normaText
$ws.Range("C21").Value = $normaText

# Row heights (ht/customHeight) per the target layout.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
